$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scalar cell fixes -----------------------------------------------
# detect_structure flag for the 2017-07-03 week goes back to "not detected"
$ws.Range("Q55").Value = 0

# isPivot flag for the 2024-06-03 week becomes a 2-bar pivot
$ws.Range("O416").Value = 2

# backup column was an empty placeholder (inline string) for the last two
# existing rows; it becomes an explicit numeric 0 once real data resumes
# below it.
$ws.Range("R418").Value = 0
$ws.Range("R419").Value = 0

# --- Append 10 new weekly rows (420-429) ------------------------------
$dateFormat = $ws.Range("A419").NumberFormat

$newRows = @(
    @{ Row=420; A=45474; B=1604.900024414062; C=1858.349975585938; D=1596.099975585938; E=1699.400024414062; F=1682.31787109375;  G=11589861; H=2024; I=7; J=1;  N=27; O=0; Q=0 },
    @{ Row=421; A=45481; B=1700.25;            C=1762.25;           D=1641.550048828125; E=1696.5;            F=1679.446899414062; G=5506325;  H=2024; I=7; J=8;  N=28; O=0; Q=0 },
    @{ Row=422; A=45488; B=1696.050048828125;  C=1777.800048828125; D=1683.550048828125; E=1698;              F=1680.931884765625; G=2428055;  H=2024; I=7; J=15; N=29; O=0; Q=1 },
    @{ Row=423; A=45495; B=1683;               C=1853.300048828125; D=1639.650024414062; E=1846.099975585938; F=1827.543212890625; G=4140171;  H=2024; I=7; J=22; N=30; O=0; Q=0 },
    @{ Row=424; A=45502; B=1850;               C=1912.150024414062; D=1805;              E=1845.849975585938; F=1827.295654296875; G=2545743;  H=2024; I=7; J=29; N=31; O=1; Q=0 },
    @{ Row=425; A=45509; B=1805;               C=1859.949951171875; D=1758.099975585938; E=1805.949951171875; F=1787.796752929688; G=2765187;  H=2024; I=8; J=5;  N=32; O=0; Q=0 },
    @{ Row=426; A=45516; B=1806;               C=1846.550048828125; D=1719.150024414062; E=1817.900024414062; F=1799.626708984375; G=1432843;  H=2024; I=8; J=12; N=33; O=0; Q=0 },
    @{ Row=427; A=45523; B=1832;               C=1874.949951171875; D=1762.5;             E=1767;              F=1767;              G=1423026;  H=2024; I=8; J=19; N=34; O=0; Q=0 },
    @{ Row=428; A=45530; B=1779;               C=1838;               D=1753.5;             E=1792.099975585938; F=1792.099975585938; G=1689448;  H=2024; I=8; J=26; N=35; O=0; Q=0 },
    @{ Row=429; A=45537; B=1794.449951171875;  C=1912;               D=1793.25;            E=1839;              F=1839;              G=2806279;  H=2024; I=9; J=2;  N=36; O=0; Q=0 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = $dateFormat
    $cellA.Value = $r.A

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = 0
    $ws.Cells.Item($row, 17).Value = $r.Q
    # Column R (backup) stays an empty placeholder for these new rows, same
    # as it was for rows 418/419 before this commit.
}
